# Bold the text in the last row of the summary tables on the closeout
# slides (Business Value Delivered, Issues & Resolutions, Enhancement
# Opportunities, Project Scope Delivered, Feature Adoption) - these
# rows were left un-bolded (only sz="1100") while the header row above
# them is bold; this brings the final table row in line intentionally.

$p = $ppt.ActivePresentation

# Slide numbers whose last-table-row-text needs to become bold.
$slideNumbers = @(4, 8, 12, 17, 21)

foreach ($slideNum in $slideNumbers) {
    $slide = $p.Slides.Item($slideNum)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $table = $shape.Table
            $lastRow = $table.Rows.Count
            for ($col = 1; $col -le $table.Columns.Count; $col++) {
                $cell = $table.Cell($lastRow, $col)
                $cell.Shape.TextFrame.TextRange.Font.Bold = 1
            }
        }
    }
}
